$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Basic Game rubric")
$ws2 = $wb.Worksheets.Item("Game extras")

# --- Sheet "Basic Game rubric" updates ---
# Row 3 - Animations: score 1 -> 2, comment updated
$ws1.Range("B3").Value = 2
$ws1.Range("C3").Value = "Mario moving & attacking, Enemies walking"

# Row 4 - Interactions: score 2 -> 3, comment updated
$ws1.Range("B4").Value = 3
$ws1.Range("C4").Value = "Player can be hit by enemies & projectiles. Player can jump on enemies, kill enemies with Pixl. Player can pickup items"

# Row 5 - Game implementation: comment updated (score unchanged)
$ws1.Range("C5").Value = "Player can move around, enemies walk around and attack, player can change attack, use items and advance to next level"

# --- Sheet "Game extras" updates ---
# Row 8 - Inventory system.: score 0 -> 1, add comment
$ws2.Range("B8").Value = 1
$ws2.Range("C8").Value = "Items and pixls"

# --- Selections / active sheet ---
$ws1.Range("B6").Select()
$ws2.Range("C9").Select()
$ws2.Activate()
